$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.381.08"
$ws.Range("E2").Value = "  +0.48%  "
$ws.Range("D3").Value = "1.876.52"
$ws.Range("E3").Value = "  +0.73%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9999"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7117"
$ws.Range("E5").Value = "  -0.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.16"
$ws.Range("E6").Value = "  +0.56%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3118"
$ws.Range("E8").Value = "  +1.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07790"
$ws.Range("E9").Value = "  +0.71%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "25.17"
$ws.Range("E10").Value = "  +0.76%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08450"
$ws.Range("E11").Value = "  +1.39%  "
$ws.Range("D12").Value = "1.872.58"
$ws.Range("E12").Value = "  +0.18%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.234"
$ws.Range("E13").Value = "  +0.57%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.7137"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.19"
$ws.Range("E15").Value = "  +0.06%  "
$ws.Range("D16").Value = "29.387.14"
$ws.Range("E16").Value = "  +0.48%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.052"
$ws.Range("E17").Value = "  +1.71%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008232"
$ws.Range("E18").Value = "  +5.09%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "241.13"
$ws.Range("E19").Value = "  -0.66%  "
$ws.Range("E20").Value = "  +0.64%  "
$ws.Range("D21").Value = "2.119.68"
$ws.Range("E21").Value = "  -0.29%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9998"
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.786"
$ws.Range("E23").Value = "  -1.40%  "
$ws.Range("E24").Value = "  +0.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1594"
$ws.Range("E25").Value = "  -0.26%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.30"
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.074"
$ws.Range("E27").Value = "  +1.96%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.54"
$ws.Range("E28").Value = "  +0.27%  "
$ws.Range("E29").Value = "  +1.02%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.425"
$ws.Range("E30").Value = "  +0.11%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.330"
$ws.Range("E31").Value = "  +2.03%  "
$ws.Range("E32").Value = "  -3.88%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05292"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.943"
$ws.Range("E34").Value = "  +0.52%  "
$ws.Range("E35").Value = "  +0.80%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7451"
$ws.Range("E36").Value = "  -9.38%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.693"
$ws.Range("E37").Value = "  +0.55%  "
$ws.Range("D39").Value = "1.231.66"
$ws.Range("E39").Value = "  +4.92%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.727"
$ws.Range("E40").Value = "  +1.23%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.479"
$ws.Range("E41").Value = "  +4.55%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "72.69"
$ws.Range("E44").Value = "  -0.20%  "
$ws.Range("E45").Value = "  +0.09%  "
$ws.Range("D46").Value = "2.017.19"
$ws.Range("E46").Value = "  -0.30%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.817"
$ws.Range("E47").Value = "  +1.51%  "
$ws.Range("E48").Value = "  +0.12%  "
$ws.Range("E49").Value = "  +2.38%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.406"
$ws.Range("E50").Value = "  +0.66%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4326"
$ws.Range("E51").Value = "  +1.29%  "

# Row 42/43 swap: Quant and TrustWalletToken exchange rows with slightly updated values
$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "110.57"
$ws.Range("E42").Value = "  +8.15%  "

$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8932"
$ws.Range("E43").Value = "  -0.49%  "
